$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1728777.2
$ws.Range("J17").Value = 1728777.2
$ws.Range("L17").Value = 5186331.6
$ws.Range("N17").Value = -5186667.6
$ws.Range("H43").Value = 6588578
$ws.Range("J43").Value = 13430
$ws.Range("L43").Value = 13430
$ws.Range("N43").Value = -13568
$ws.Range("H116").Value = 2700
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents() | Out-Null
$ws.Range("H137").Value = 1878.7812
$ws.Range("I137").Value = 1878.7812
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 5636.3436
$ws.Range("L137").Value = 0
$ws.Range("M137").ClearContents() | Out-Null
$ws.Range("N137").Value = -3086.3436
$ws.Range("H138").Value = 424627.97
$ws.Range("J138").Value = 911125.4399999999
$ws.Range("L138").Value = 2733376.32
$ws.Range("N138").Value = -2743656.32

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 503.6
$ws.Range("I4").Value = 732.6667
$ws.Range("J4").Value = 160
$ws.Range("K4").Value = 732.6667
$ws.Range("L4").Value = 160
$ws.Range("M4").Value = -616.6667
$ws.Range("N4").Value = -392
$ws.Range("H45").Value = 2954.2727
$ws.Range("J45").Value = 3004.6667
$ws.Range("L45").Value = 3004.6667
$ws.Range("N45").Value = -3758.6667
$ws.Range("H122").Value = 5133.2
$ws.Range("I122").Value = 4619.857
$ws.Range("K122").Value = 13859.571
$ws.Range("M122").Value = -11409.571
$ws.Range("H132").Value = 1510.5193
$ws.Range("I132").Value = 1228.4681
$ws.Range("K132").Value = 3685.4043
$ws.Range("M132").Value = -1155.4043

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 9845.799999999999
$ws.Range("I86").Value = 16348.125
$ws.Range("J86").Value = 2414.5715
$ws.Range("K86").Value = 16348.125
$ws.Range("L86").Value = 2414.5715
$ws.Range("M86").Value = -15225.125
$ws.Range("N86").Value = -4660.5715
$ws.Range("H89").Value = 9845.799999999999
$ws.Range("I89").Value = 16348.125
$ws.Range("J89").Value = 2414.5715
$ws.Range("K89").Value = 81740.625
$ws.Range("L89").Value = 12072.8575
$ws.Range("M89").Value = -76124.625
$ws.Range("N89").Value = -23304.8575
$ws.Range("H99").Value = 2932.9333
$ws.Range("I99").Value = 2440.5293
$ws.Range("J99").Value = 3576.8462
$ws.Range("K99").Value = 2440.5293
$ws.Range("L99").Value = 3576.8462
$ws.Range("M99").Value = -942.5293000000001
$ws.Range("N99").Value = -6572.8462
$ws.Range("H105").Value = 3417.5862
$ws.Range("I105").Value = 2049.9583
$ws.Range("K105").Value = 2049.9583
$ws.Range("M105").Value = -302.9582999999998
$ws.Range("H134").Value = 1777.7317
$ws.Range("I134").Value = 1098.3636
$ws.Range("K134").Value = 3295.0908
$ws.Range("M134").Value = -760.0907999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1967.2972
$ws.Range("I31").Value = 1266.4242
$ws.Range("J31").Value = 7749.5
$ws.Range("K31").Value = 1266.4242
$ws.Range("L31").Value = 7749.5
$ws.Range("M31").Value = -971.4241999999999
$ws.Range("N31").Value = -8339.5
$ws.Range("H34").Value = 1967.2972
$ws.Range("I34").Value = 1266.4242
$ws.Range("J34").Value = 7749.5
$ws.Range("K34").Value = 1266.4242
$ws.Range("L34").Value = 7749.5
$ws.Range("M34").Value = -1064.4242
$ws.Range("N34").Value = -8153.5
$ws.Range("H58").Value = 1687.7778
$ws.Range("I58").Value = 1733
$ws.Range("K58").Value = 1733
$ws.Range("M58").Value = -1530
$ws.Range("H96").Value = 12725
$ws.Range("J96").Value = 12725
$ws.Range("L96").Value = 12725
$ws.Range("N96").Value = -18217
$ws.Range("H107").Value = 1699.625
$ws.Range("I107").Value = 696.5
$ws.Range("J107").Value = 2034
$ws.Range("K107").Value = 696.5
$ws.Range("L107").Value = 2034
$ws.Range("M107").Value = 1223.5
$ws.Range("N107").Value = -5874
$ws.Range("H135").Value = 88500
$ws.Range("J135").Value = 83333.336
$ws.Range("L135").Value = 83333.336
$ws.Range("N135").Value = -93473.336
$ws.Range("H136").Value = 1687.7778
$ws.Range("I136").Value = 1733
$ws.Range("K136").Value = 5199
$ws.Range("M136").Value = -2649

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 655.93335
$ws.Range("I12").Value = 1317.5
$ws.Range("J12").Value = 415.36365
$ws.Range("K12").Value = 3952.5
$ws.Range("L12").Value = 1246.09095
$ws.Range("M12").Value = -3779.5
$ws.Range("N12").Value = -1592.09095
$ws.Range("H33").Value = 350.46667
$ws.Range("I33").Value = 184.33333
$ws.Range("J33").Value = 461.22223
$ws.Range("K33").Value = 1105.99998
$ws.Range("L33").Value = 2767.33338
$ws.Range("M33").Value = -822.9999800000001
$ws.Range("N33").Value = -3333.33338
$ws.Range("H68").Value = 41751
$ws.Range("I68").Value = 72736
$ws.Range("J68").Value = 4569
$ws.Range("K68").Value = 218208
$ws.Range("L68").Value = 13707
$ws.Range("M68").Value = -217397
$ws.Range("N68").Value = -15329
$ws.Range("H71").Value = 41751
$ws.Range("I71").Value = 72736
$ws.Range("J71").Value = 4569
$ws.Range("K71").Value = 654624
$ws.Range("L71").Value = 41121
$ws.Range("M71").Value = -650568
$ws.Range("N71").Value = -49233
$ws.Range("H139").Value = 3268.68
$ws.Range("I139").Value = 2033.1428
$ws.Range("K139").Value = 6099.428400000001
$ws.Range("M139").Value = -959.4284000000007

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2501.2273
$ws.Range("I22").Value = 932.8570999999999
$ws.Range("K22").Value = 932.8570999999999
$ws.Range("M22").Value = -637.8570999999999
$ws.Range("H27").Value = 2501.2273
$ws.Range("I27").Value = 932.8570999999999
$ws.Range("K27").Value = 932.8570999999999
$ws.Range("M27").Value = -825.8570999999999
$ws.Range("H46").Value = 8169.9473
$ws.Range("J46").Value = 11975.363
$ws.Range("L46").Value = 11975.363
$ws.Range("N46").Value = -12351.363
$ws.Range("H55").Value = 209.86363
$ws.Range("I55").Value = 183.07692
$ws.Range("J55").Value = 248.55556
$ws.Range("K55").Value = 183.07692
$ws.Range("L55").Value = 248.55556
$ws.Range("M55").Value = -10.07692
$ws.Range("N55").Value = -594.55556
$ws.Range("H132").Value = 3113.2927
$ws.Range("I132").Value = 3088.4849
$ws.Range("K132").Value = 9265.4547
$ws.Range("M132").Value = -6735.4547
$ws.Range("H136").Value = 3950.1333
$ws.Range("I136").Value = 3667.25
$ws.Range("K136").Value = 11001.75
$ws.Range("M136").Value = -8451.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 2634.5881
$ws.Range("I100").Value = 2253.9
$ws.Range("J100").Value = 3178.4285
$ws.Range("K100").Value = 4507.8
$ws.Range("L100").Value = 6356.857
$ws.Range("M100").Value = -3966.8
$ws.Range("N100").Value = -7438.857
$ws.Range("H113").Value = 2688882.2
$ws.Range("J113").Value = 1174.8572
$ws.Range("L113").Value = 3524.5716
$ws.Range("N113").Value = -7864.571599999999
$ws.Range("H132").Value = 2034.5172
$ws.Range("I132").Value = 1982.6531
$ws.Range("J132").Value = 2316.889
$ws.Range("K132").Value = 5947.9593
$ws.Range("L132").Value = 6950.667
$ws.Range("M132").Value = -3417.9593
$ws.Range("N132").Value = -12010.667
